# "added test for allocrule2 fix"
# Appends a new ftest case row (fm43) to the "ftests" sheet, describing a
# test for multiple accounts with a differing number of layers (policies)
# per account at account-level output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Restore the window/view position & size the author's Excel instance had
# on save, and scroll so row 21 is pinned at the top (best-effort - some
# hosts don't round-trip window geometry, but it's harmless to set).
$win = $wb.Windows.Item(1)
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12576
$win.ScrollRow = 21
$win.ScrollColumn = 1

# The existing test rows run from row 4 to row 47 (B:I, sometimes B:K).
# Row 47 is the last populated row and carries the per-column formatting
# (borders/shading/number styles) that every data row re-uses, so copy its
# formatting down into the new row 48 before filling in the new values -
# this keeps the same style index per column as every other data row.
$lastRow = 47
$newRow = 48

$ws.Range("B" + $lastRow + ":I" + $lastRow).Copy() | Out-Null
$ws.Range("B" + $newRow + ":I" + $newRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 2).Value = "fm43"
$ws.Cells.Item($newRow, 3).Value = "Multiple accounts with different number of layers (policies) per account. Account level output"
$ws.Cells.Item($newRow, 4).Value = "All"
$ws.Cells.Item($newRow, 5).Value = "2,12,14"
$ws.Cells.Item($newRow, 6).Value = 3
$ws.Cells.Item($newRow, 7).Value = "1,2"
$ws.Cells.Item($newRow, 8).Value = "in progress"
$ws.Cells.Item($newRow, 9).Value = "in progress"

# Match the author's on-save selection/view state.
$ws.Activate() | Out-Null
$ws.Range("E49").Select() | Out-Null

Write-Host "Added fm43 test case row to 'ftests' sheet"
